$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: paragraph "Väl inne på schemasidan ... framtida aktivitetstillfällen.
# Om du är registrerad som ett par i schemat kan du ange 1,2, Nej eller Kanske,
# annars Ja, Nej eller Kanske."
# becomes two paragraphs, with extra explanatory text and a _GoBack bookmark
# around "Alternativet "Kanske" finns också med för att ange att du vill
# bestämma dig senare".
# -----------------------------------------------------------------------

$dq = [char]0x201D

$oldTail = " aktivitetstillfällen. Om du är registrerad som ett par i schemat kan du ange 1,2, Nej eller Kanske, annars Ja, Nej eller Kanske."
$newP1Tail = " aktivitetstillfällen. Om du är registrerad som ett par i schemat kan du ange om ingen, en eller två personer kommer."
$newP2 = " Annars registrerar du Ja eller Nej. Alternativet " + $dq + "Kanske" + $dq + " finns också med för att ange att du vill bestämma dig senare."

$rng = $d.Content
$rng.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = $newP1Tail + [char]13 + $newP2

# Split the tail of paragraph 3 into three runs: "... ange ", "om ingen, en
# eller " and "två personer kommer."
$p3 = $d.Paragraphs(3)
$splitA = $p3.Range.Duplicate
$splitA.Find.Execute("om ingen, en eller två personer kommer.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("ZZtmp1", $splitA)
$splitB = $d.Bookmarks("ZZtmp1").Range.Duplicate
$splitB.Find.Execute("om ingen, en eller ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("ZZtmp2", $splitB)
$d.Bookmarks("ZZtmp1").Delete()
$d.Bookmarks("ZZtmp2").Delete()

# The new paragraph 4 currently holds one run with all of $newP2's text
# (minus the leading paragraph mark). Give it the run layout from the diff:
# " " / "Annars registrerar du Ja eller Nej. " / (bookmarked) "Alternativet
# ...senare" / "."
$p4 = $d.Paragraphs(4)
$altText = "Alternativet " + $dq + "Kanske" + $dq + " finns också med för att ange att du vill bestämma dig senare"
$altRng = $p4.Range.Duplicate
$altRng.Find.Execute($altText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $altRng)

$p4b = $d.Paragraphs(4)
$leadRng = $p4b.Range.Duplicate
$leadRng.Find.Execute("Annars registrerar du Ja eller Nej. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("ZZtmp3", $leadRng)
$d.Bookmarks("ZZtmp3").Delete()

# -----------------------------------------------------------------------
# Change 2: "systemadministratörsbehörighet" -> "schemaansvarsbehörighet",
# keeping the two existing runs ("system" / "administratörsbehörighet...")
# as two runs ("schemaansva" / "rsbehörighet...").
# (Paragraph index is 6, not 5, because of the paragraph inserted above.)
# -----------------------------------------------------------------------

$p5 = $d.Paragraphs(6)
$sysRng = $p5.Range.Duplicate
$sysRng.Find.Execute("system", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("ZZtmp4", $sysRng)
$d.Bookmarks("ZZtmp4").Range.Text = "schemaansva"

$p5b = $d.Paragraphs(6)
$adminRng = $p5b.Range.Duplicate
$adminRng.Find.Execute("administratörsbehörighet", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$adminRng.Text = "rsbehörighet"
$d.Bookmarks("ZZtmp4").Delete()

# -----------------------------------------------------------------------
# Change 3: "Systemadministratör kan " -> "Schemaansvarig kan ", keeping the
# existing three runs ("En " / "System" / "administratör kan ") as three runs
# ("En " / "Schemaansvarig" / " kan ").
# (Paragraph index is 7, not 6, because of the paragraph inserted above.)
# -----------------------------------------------------------------------

$p6 = $d.Paragraphs(7)
$sys2Rng = $p6.Range.Duplicate
$sys2Rng.Find.Execute("System", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("ZZtmp5", $sys2Rng)
$d.Bookmarks("ZZtmp5").Range.Text = "Schemaansvarig"

$p6b = $d.Paragraphs(7)
$admin2Rng = $p6b.Range.Duplicate
$admin2Rng.Find.Execute("administratör kan ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$admin2Rng.Text = " kan "
$d.Bookmarks("ZZtmp5").Delete()

# -----------------------------------------------------------------------
# Change 4: merge the three runs "a" / "nge om medlemmen representerar en
# eller två personer" / " i respektive schema" into a single run.
# -----------------------------------------------------------------------

$mergedText = "ange om medlemmen representerar en eller två personer i respektive schema"
$rngMerge = $d.Content
$rngMerge.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)

# -----------------------------------------------------------------------
# Change 5: remove the old (now-redundant) _GoBack bookmark that used to sit
# after "registrera medlemmar (bör dock undvikas)".
# -----------------------------------------------------------------------

$oldGoBack = $d.Range(0, 0)
$found = $false
$searchRng = $d.Content
$searchRng.Find.Execute("registrera medlemmar (bör dock undvikas)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterRng = $d.Range($searchRng.End, $searchRng.End)

if ($d.Bookmarks.Exists("_GoBack")) {
    $gb = $d.Bookmarks("_GoBack")
    if ($gb.Start -eq $afterRng.Start) {
        $gb.Delete()
    }
}

Write-Output "edit complete"
